$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the two runs that spell "...konkretne style – w" + "tedy zmieni..."
#    back into a single run "...konkretne style – wtedy zmieni...". The
#    original had a (now stale) "_GoBack" bookmark splitting the two runs;
#    replacing the whole (contiguous) text in one go merges the runs and
#    drops that bookmark automatically.
# ---------------------------------------------------------------------------
$mergedText = "Aby modyfikować style poszczególnych sekcji, najlepiej robić to modyfikując konkretne style – wtedy zmieni się styl dla np. wszystkich akordów lub wszystkich tytułów."
$null = $d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)

# ---------------------------------------------------------------------------
# 2. Re-create the "_GoBack" bookmark right after the word "Tonacja" (this is
#    where the cursor was left after the last edit). Adding a bookmark to a
#    zero-length Range positioned exactly at that paragraph's end is
#    unreliable, so nudge a temporary marker character in, anchor the
#    bookmark next to it, then remove the marker again.
# ---------------------------------------------------------------------------
$tonacja = $d.Content
$null = $tonacja.Find.Execute("Tonacja", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tonacja.Collapse(0)
$tonacja.InsertAfter("\x01")
$tonacja.Collapse(1)
$d.Bookmarks.Add("_GoBack", $tonacja)

$marker = $d.Content
$null = $marker.Find.Execute("\x01", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$marker.Delete()

# ---------------------------------------------------------------------------
# 3. "Autor piosenki" / "Autor piosenki Znak" styles become italic.
# ---------------------------------------------------------------------------
$d.Styles("Autor piosenki").Font.Italic = $true
$d.Styles("Autor piosenki Znak").Font.Italic = $true
